$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 150.5
$ws.Range("I12").Value = 150.5
$ws.Range("K12").Value = 150.5
$ws.Range("M12").Value = 19.5

$ws.Range("H40").Value = 4603.778
$ws.Range("I40").Value = 3405.6667
$ws.Range("J40").Value = 7000
$ws.Range("K40").Value = 3405.6667
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = -3230.6667
$ws.Range("N40").Value = -7350

$ws.Range("H41").Value = 219.1
$ws.Range("I41").Value = 164.57143
$ws.Range("K41").Value = 164.57143
$ws.Range("M41").Value = 275.42857

$ws.Range("H80").Value = 2706.111
$ws.Range("I80").Value = 2979.5715
$ws.Range("K80").Value = 8938.7145
$ws.Range("M80").Value = -7940.7145

$ws.Range("H83").Value = 2706.111
$ws.Range("I83").Value = 2979.5715
$ws.Range("K83").Value = 26816.1435
$ws.Range("M83").Value = -21824.1435

$ws.Range("H100").Value = 2799.8
$ws.Range("I100").Value = 10000
$ws.Range("J100").Value = 999.75
$ws.Range("K100").Value = 10000
$ws.Range("L100").Value = 999.75
$ws.Range("M100").Value = -9459
$ws.Range("N100").Value = -2081.75

$ws.Range("H101").Value = 2099
$ws.Range("J101").Value = 499.33334
$ws.Range("L101").Value = 1498.00002
$ws.Range("N101").Value = -4742.000019999999

$ws.Range("H111").Value = 10032
$ws.Range("J111").Value = 10032
$ws.Range("L111").Value = 30096
$ws.Range("N111").Value = -36230

$ws.Range("H137").Value = 1725
$ws.Range("I137").Value = 1725
$ws.Range("K137").Value = 5175
$ws.Range("M137").Value = -2625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1579.8
$ws.Range("I61").Value = 2450
$ws.Range("J61").Value = 999.6667
$ws.Range("K61").Value = 2450
$ws.Range("L61").Value = 999.6667
$ws.Range("M61").Value = -2238
$ws.Range("N61").Value = -1423.6667

$ws.Range("H110").Value = 2193.077
$ws.Range("I110").Value = 2237
$ws.Range("J110").Value = 1951.5
$ws.Range("K110").Value = 2237
$ws.Range("L110").Value = 1951.5
$ws.Range("M110").Value = -192
$ws.Range("N110").Value = -6041.5

$ws.Range("H132").Value = 3517.75
$ws.Range("I132").Value = 3024
$ws.Range("K132").Value = 9072
$ws.Range("M132").Value = -6542

$ws.Range("H136").Value = 1579.8
$ws.Range("I136").Value = 2450
$ws.Range("J136").Value = 999.6667
$ws.Range("K136").Value = 7350
$ws.Range("L136").Value = 2999.0001
$ws.Range("M136").Value = -4800
$ws.Range("N136").Value = -8099.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5139.2
$ws.Range("I20").Value = 1333
$ws.Range("K20").Value = 1333
$ws.Range("M20").Value = -1086

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H94").Value = 7995
$ws.Range("I94").Value = 990
$ws.Range("J94").Value = 15000
$ws.Range("K94").Value = 990
$ws.Range("L94").Value = 15000
$ws.Range("M94").Value = -539
$ws.Range("N94").Value = -15902

$ws.Range("H105").Value = 6084.6665
$ws.Range("J105").Value = 3500
$ws.Range("L105").Value = 3500
$ws.Range("N105").Value = -6994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2825
$ws.Range("J31").Value = 2931.3333
$ws.Range("L31").Value = 2931.3333
$ws.Range("N31").Value = -3521.3333

$ws.Range("H34").Value = 2825
$ws.Range("J34").Value = 2931.3333
$ws.Range("L34").Value = 2931.3333
$ws.Range("N34").Value = -3335.3333

$ws.Range("H107").Value = 1720.6428
$ws.Range("I107").Value = 1953.6364
$ws.Range("K107").Value = 1953.6364
$ws.Range("M107").Value = -33.63640000000009

$ws.Range("H132").Value = 3449.5
$ws.Range("J132").Value = 3449.5
$ws.Range("L132").Value = 10348.5
$ws.Range("N132").Value = -15408.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 620.6667
$ws.Range("I97").Value = 569.875
$ws.Range("K97").Value = 1709.625
$ws.Range("M97").Value = -1213.625

$ws.Range("H113").Value = 1374.3334
$ws.Range("I113").Value = 623.5
$ws.Range("K113").Value = 1870.5
$ws.Range("M113").Value = 299.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 397.14285
$ws.Range("I16").Value = 398
$ws.Range("J16").Value = 395
$ws.Range("K16").Value = 398
$ws.Range("L16").Value = 395
$ws.Range("M16").Value = -228
$ws.Range("N16").Value = -735

$ws.Range("H122").Value = 5121.3125
$ws.Range("I122").Value = 5121.3125
$ws.Range("K122").Value = 15363.9375
$ws.Range("M122").Value = -12913.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 676666.7
$ws.Range("I2").Value = 2000000
$ws.Range("K2").Value = 2000000
$ws.Range("M2").Value = -1999888

$ws.Range("H4").Value = 20000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H41").Value = 23294.666
$ws.Range("I41").Value = 28943
$ws.Range("K41").Value = 28943
$ws.Range("M41").Value = -28553

$ws.Range("H122").Value = 3036.4666
$ws.Range("I122").Value = 2982.0715
$ws.Range("J122").Value = 3798
$ws.Range("K122").Value = 8946.2145
$ws.Range("L122").Value = 11394
$ws.Range("M122").Value = -6496.2145
$ws.Range("N122").Value = -16294

$ws.Range("H132").Value = 2186.5
$ws.Range("I132").Value = 1691.7142
$ws.Range("J132").Value = 2879.2
$ws.Range("K132").Value = 5075.142599999999
$ws.Range("L132").Value = 8637.599999999999
$ws.Range("M132").Value = -2545.142599999999
$ws.Range("N132").Value = -13697.6
